# Weekly fruit/vegetable price update.
# Insert two new daily records at the top of the historical data block
# (rows 30:31), pushing all subsequent rows down by two. This naturally
# relocates the former rows 51/52 to rows 53/54, matching the target
# state without any further edits needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("30:31").Insert()

# New row 30
$ws.Cells.Item(30, 1).Value = 8
$ws.Cells.Item(30, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44566
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = 100112030
$ws.Cells.Item(30, 7).Value = "Poroto granado"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 400
$ws.Cells.Item(30, 11).Value = 16000
$ws.Cells.Item(30, 12).Value = 17000
$ws.Cells.Item(30, 13).Value = 16500
$ws.Cells.Item(30, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(30, 16).Value = 1100
$ws.Cells.Item(30, 17).Value = 15
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# New row 31
$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44566
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112030
$ws.Cells.Item(31, 7).Value = "Poroto granado"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 600
$ws.Cells.Item(31, 11).Value = 31000
$ws.Cells.Item(31, 12).Value = 32000
$ws.Cells.Item(31, 13).Value = 31500
$ws.Cells.Item(31, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(31, 16).Value = 1260
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
